$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1: add a new "Checklista innan handling" feature row, move the
#    selection, and widen column C to fit the widest entry.
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("C4").Value = "Checklista innan handling"
$sheet1.Range("D4").Value = "Möjlighet att kunna få en checklista på de varor man har I skafferiet osv. "
$sheet1.Columns.Item(3).ColumnWidth = 58.5
$sheet1.Range("D5").Select()

# ---------------------------------------------------------------------------
# 2. Databas sheet: rework the "recipe" data-model documentation block so
#    ingredients become their own top-level entity referenced from a new
#    "subRecipe" entity.
# ---------------------------------------------------------------------------
$databas = $wb.Worksheets.Item("Databas")

# K2 used to be "ingridient" (child of recipe) -> becomes "subRecipe"
$databas.Range("K2").Value = "subRecipe"

# The ingredient entity moves out to its own column (O) with its own fields
$databas.Range("O2").Value = "ingridient"
$databas.Range("O3").Value = "id"
$databas.Range("O4").Value = "category"
$databas.Range("O5").Value = "name"

# "category" entity moves from column M to column Q
$databas.Range("Q2").Value = "category"
$databas.Range("Q3").Value = "id"
$databas.Range("M2").ClearContents()
$databas.Range("M3").ClearContents()
$databas.Range("K3").ClearContents()
$databas.Range("K4").ClearContents()

# recipe's "list of ingridents" field becomes "subrecipes", and gains
# difficulty / name / portions fields
$databas.Range("I4").Value = "subrecipes"
$databas.Range("I6").Value = "difficulty"
$databas.Range("I7").Value = "name"
$databas.Range("I8").Value = "portions"

$databas.Range("I2:I8").Select()
$databas.Activate()

# ---------------------------------------------------------------------------
# 3. New "TestRecept" sheet: a real recipe row backed by a table, plus a
#    column of the recipe field names used above.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testRecept = $wb.Worksheets.Add($null, $lastSheet)
$testRecept.Name = "TestRecept"

$testRecept.Range("F1").Value = "Name"
$testRecept.Range("G1").Value = "Time h"
$testRecept.Range("H1").Value = "NumberOfPortions"
$testRecept.Range("I1").Value = "Difficulty"
$testRecept.Range("J1").Value = "Instructions"
$testRecept.Range("K1").Value = "Image link"

$testRecept.Range("F2").Value = "Färgglad pastasallad med fläskfile"
$testRecept.Range("G2").Value = 0.5
$testRecept.Range("H2").Value = 4
$testRecept.Range("I2").Value = 2
$testRecept.Range("J2").Value = "Putsa köttet från senor och hinnor. Skär det i 2 cm tjocka bitar. Blanda vetemjöl, salt och peppar och vänd köttet i det.`nBryn köttet på båda sidor i en medelvarm panna. Stek dem därefter 2-3 minuter tills köttet är genomstekt. `nKoka pastan enligt anvisning på paketet, men utan salt, eftersom rätten i övrigt ger så mycket smak. Låt sockerärter eller haricotsverts koka med sista 2-3 minuterna. Häll av i durkslag och skölj med kallt vatten. `nRiv morötterna. Skär gurkan i tärningar och tomaterna i bitar. Blanda alla grönsaker med pastan.`nBlanda dressingen och häll över. Fördela köttet över salladen och servera med bröd."
$testRecept.Range("J2").WrapText = $true
$testRecept.Rows.Item(2).RowHeight = 43.5

$testRecept.Range("C3").Value = "recipe"
$testRecept.Range("C4").Value = "time"
$testRecept.Range("C5").Value = "subrecipes"
$testRecept.Range("C6").Value = "instructions"
$testRecept.Range("C7").Value = "difficulty"
$testRecept.Range("C8").Value = "name"
$testRecept.Range("C9").Value = "portions"

$lo = $testRecept.ListObjects.Add(1, $testRecept.Range("F1:K2"), $null, 1)
$lo.Name = "Table1"

$testRecept.Columns.Item(6).ColumnWidth = 31.42578125
$testRecept.Columns.Item(8).ColumnWidth = 20
$testRecept.Columns.Item(9).ColumnWidth = 11.28515625
$testRecept.Columns.Item(10).ColumnWidth = 13.5703125

$testRecept.Range("M14").Select()

# Leave "Databas" as the active tab, matching the authored workbook state.
$databas.Activate()
